$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 3 - shifts the old data rows (3..12) down
#    to (4..13) and pushes the trailing spacer row from 23 to 24.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").Insert()

# ---------------------------------------------------------------------------
# 2. Header row (row 1) - rename / re-purpose a few headers.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "LHV"
$ws.Range("L1").Value = "biomass-to-fuel ratio"
$ws.Range("M1").Value = "CO2 removal"
$ws.Range("N1").Value = "is fossil"
$ws.Range("O1").Value = "is biofuel"
$ws.Range("P1").Value = "meta-notes"
$ws.Range("Q1").Value = "meta-source"

# ---------------------------------------------------------------------------
# 3. Units row (row 2) - update text, drop two stale unit cells, add new one,
#    and apply the new italic "meta" style to every populated cell.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "(gj/t dry)"
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = "t CO2 / t fuel source"
$ws.Range("N2").ClearContents()

$ws.Range("A2").Font.Italic = $true
$ws.Range("B2").Font.Italic = $true
$ws.Range("C2").Font.Italic = $true
$ws.Range("D2").Font.Italic = $true
$ws.Range("E2").Font.Italic = $true
$ws.Range("F2").Font.Italic = $true
$ws.Range("K2").Font.Italic = $true
$ws.Range("M2").Font.Italic = $true

# ---------------------------------------------------------------------------
# 4. New row 3 - a blank "meta-notes" template row.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "meta-notes"
$ws.Range("M3").Value = 0
$ws.Range("A3").Font.Italic = $true
$ws.Range("M3").Font.Italic = $true

# ---------------------------------------------------------------------------
# 5. Data rows (now 4..13) - add the new lookup/flag columns, fix up the
#    upstream-CO2 formulas for natural gas & diesel, rename the Eurofer row,
#    and relocate the old "meta-source" tags from column P to the new
#    column Q.
# ---------------------------------------------------------------------------

# row 4 - coal
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 1
$ws.Range("O4").Formula = "=1-N4"
$ws.Range("P4").ClearContents()
$ws.Range("Q4").Value = "ecoinvent 2.2"

# row 5 - charcoal
$ws.Range("M5").Value = 1.25
$ws.Range("N5").Value = 0
$ws.Range("M5").NumberFormat = "0.00"
$ws.Range("N5").NumberFormat = "0.00"
$ws.Range("O5").Formula = "=1-N5"
$ws.Range("P5").ClearContents()
$ws.Range("Q5").Value = "ecoinvent 2.2"

# row 6 - natural gas
$ws.Range("K6").Formula = "=0.007*C6"
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 1
$ws.Range("O6").Formula = "=1-N6"
$ws.Range("Q6").Value = "ecoinvent 2.2"

# row 7 - diesel
$ws.Range("K7").Value = 0.438
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 1
$ws.Range("O7").Formula = "=1-N7"
$ws.Range("Q7").Value = "ecoinvent 2.2"

# row 8 - wood chips
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 1.25
$ws.Range("N8").Value = 0
$ws.Range("M8").NumberFormat = "0.00"
$ws.Range("N8").NumberFormat = "0.00"
$ws.Range("O8").Formula = "=1-N8"

# row 9 - steam
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 1
$ws.Range("O9").Formula = "=1-N9"

# row 10 - electricity-eurofer (renamed from "Eurofer electricity mix proxy")
$ws.Range("A10").Value = "electricity-eurofer"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 1
$ws.Range("O10").Formula = "=1-N10"
$ws.Range("P10").ClearContents()
$ws.Range("Q10").Value = "EUROFER"

# row 11 - coke
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 1
$ws.Range("O11").Formula = "=1-N11"
$ws.Range("P11").ClearContents()
$ws.Range("Q11").Value = "IEAGHG 2013"

# row 12 - coking coal
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 1
$ws.Range("O12").Formula = "=1-N12"
$ws.Range("P12").ClearContents()
$ws.Range("Q12").Value = "IEAGHG 2013"

# row 13 - PCI coal
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 1
$ws.Range("O13").Formula = "=1-N13"
$ws.Range("P13").ClearContents()
$ws.Range("Q13").Value = "IEAGHG"

# ---------------------------------------------------------------------------
# 6. Defined name - the "fuels" table now spans one extra row.
# ---------------------------------------------------------------------------
$wb.Names.Item("Fuels!fuels").RefersTo = "=Fuels!`$A`$1:`$E`$8"
